$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Insert 16 new rows starting at row 153 (shifts old rows 153-155 down to 169-171)
$ws.Range("A153:A168").EntireRow.Insert()

    # Row 153
    $ws.Cells.Item(153, 1).Value = 152
    $ws.Cells.Item(153, 3).Value = 45252
    $ws.Cells.Item(153, 4).Value = '2023-11-22'

    # Row 154
    $ws.Cells.Item(154, 1).Value = 153
    $ws.Cells.Item(154, 3).Value = 45253
    $ws.Cells.Item(154, 4).Value = '2023-11-23'

    # Row 155
    $ws.Cells.Item(155, 1).Value = 154
    $ws.Cells.Item(155, 2).Value = 144
    $ws.Cells.Item(155, 3).Value = 45254
    $ws.Cells.Item(155, 4).Value = '2023-11-24'
    $ws.Cells.Item(155, 5).Value = 1
    $ws.Cells.Item(155, 6).Value = 1.07
    $ws.Cells.Item(155, 7).Formula = '=I152'
    $ws.Cells.Item(155, 8).Value = 330
    $ws.Cells.Item(155, 9).Formula = '=G155+H155'
    $ws.Cells.Item(155, 10).Value = 'ESPORTS'
    $ws.Cells.Item(155, 11).Value = 'CS2 ESPORTS BATTLE'
    $ws.Cells.Item(155, 12).Formula = '=ROUND((I155/$G$31-1)*100,3)+$L$29'

    # Row 156
    $ws.Cells.Item(156, 1).Value = 155
    $ws.Cells.Item(156, 2).Value = 145
    $ws.Cells.Item(156, 3).Value = 45254
    $ws.Cells.Item(156, 4).Value = '2023-11-24'
    $ws.Cells.Item(156, 5).Value = 1
    $ws.Cells.Item(156, 6).Value = 1.05
    $ws.Cells.Item(156, 7).Formula = '=I155'
    $ws.Cells.Item(156, 8).Value = 250
    $ws.Cells.Item(156, 9).Formula = '=G156+H156'
    $ws.Cells.Item(156, 10).Value = 'ESPORTS'
    $ws.Cells.Item(156, 11).Value = 'DOTA 2 ULTRAS DOTA PRO'
    $ws.Cells.Item(156, 12).Formula = '=ROUND((I156/$G$31-1)*100,3)+$L$29'

    # Row 157
    $ws.Cells.Item(157, 1).Value = 156
    $ws.Cells.Item(157, 3).Value = 45255
    $ws.Cells.Item(157, 4).Value = '2023-11-25'

    # Row 158
    $ws.Cells.Item(158, 1).Value = 157
    $ws.Cells.Item(158, 2).Value = 146
    $ws.Cells.Item(158, 3).Value = 45256
    $ws.Cells.Item(158, 4).Value = '2023-11-26'
    $ws.Cells.Item(158, 5).Value = 1
    $ws.Cells.Item(158, 6).Value = 1.168
    $ws.Cells.Item(158, 7).Formula = '=I156'
    $ws.Cells.Item(158, 8).Value = 880
    $ws.Cells.Item(158, 9).Formula = '=G158+H158'
    $ws.Cells.Item(158, 10).Value = 'BASKET'
    $ws.Cells.Item(158, 11).Value = 'NBA'
    $ws.Cells.Item(158, 12).Formula = '=ROUND((I158/$G$31-1)*100,3)+$L$29'

    # Row 159
    $ws.Cells.Item(159, 1).Value = 158
    $ws.Cells.Item(159, 3).Value = 45257
    $ws.Cells.Item(159, 4).Value = '2023-11-27'

    # Row 160
    $ws.Cells.Item(160, 1).Value = 159
    $ws.Cells.Item(160, 3).Value = 45258
    $ws.Cells.Item(160, 4).Value = '2023-11-28'

    # Row 161
    $ws.Cells.Item(161, 1).Value = 160
    $ws.Cells.Item(161, 2).Value = 147
    $ws.Cells.Item(161, 3).Value = 45259
    $ws.Cells.Item(161, 4).Value = '2023-11-29'
    $ws.Cells.Item(161, 5).Value = 1
    $ws.Cells.Item(161, 6).Value = 1.112
    $ws.Cells.Item(161, 7).Formula = '=I158'
    $ws.Cells.Item(161, 8).Value = 180
    $ws.Cells.Item(161, 9).Formula = '=G161+H161'
    $ws.Cells.Item(161, 10).Value = 'TENIS DE MESA'
    $ws.Cells.Item(161, 11).Value = 'MASTERS'
    $ws.Cells.Item(161, 12).Formula = '=ROUND((I161/$G$31-1)*100,3)+$L$29'

    # Row 162
    $ws.Cells.Item(162, 1).Value = 161
    $ws.Cells.Item(162, 2).Value = 148
    $ws.Cells.Item(162, 3).Value = 45259
    $ws.Cells.Item(162, 4).Value = '2023-11-29'
    $ws.Cells.Item(162, 5).Value = 1
    $ws.Cells.Item(162, 6).Value = 1.165
    $ws.Cells.Item(162, 7).Formula = '=I161'
    $ws.Cells.Item(162, 8).Value = 495
    $ws.Cells.Item(162, 9).Formula = '=G162+H162'
    $ws.Cells.Item(162, 10).Value = 'TENIS DE MESA'
    $ws.Cells.Item(162, 11).Value = 'MASTERS'
    $ws.Cells.Item(162, 12).Formula = '=ROUND((I162/$G$31-1)*100,3)+$L$29'

    # Row 163
    $ws.Cells.Item(163, 1).Value = 162
    $ws.Cells.Item(163, 3).Value = 45260
    $ws.Cells.Item(163, 4).Value = '2023-11-30'

    # Row 164
    $ws.Cells.Item(164, 1).Value = 163
    $ws.Cells.Item(164, 3).Value = 45261
    $ws.Cells.Item(164, 4).Value = '2023-12-01'

    # Row 165
    $ws.Cells.Item(165, 1).Value = 164
    $ws.Cells.Item(165, 3).Value = 45262
    $ws.Cells.Item(165, 4).Value = '2023-12-02'

    # Row 166
    $ws.Cells.Item(166, 1).Value = 165
    $ws.Cells.Item(166, 3).Value = 45263
    $ws.Cells.Item(166, 4).Value = '2023-12-03'

    # Row 167
    $ws.Cells.Item(167, 1).Value = 166
    $ws.Cells.Item(167, 2).Value = 149
    $ws.Cells.Item(167, 3).Value = 45264
    $ws.Cells.Item(167, 4).Value = '2023-12-04'
    $ws.Cells.Item(167, 5).Value = 0
    $ws.Cells.Item(167, 6).Value = 1.152
    $ws.Cells.Item(167, 7).Formula = '=I162'
    $ws.Cells.Item(167, 8).Value = -1887
    $ws.Cells.Item(167, 9).Formula = '=G167+H167'
    $ws.Cells.Item(167, 10).Value = 'ESPORTS'
    $ws.Cells.Item(167, 11).Value = 'CS2 CCT ONLINE FINALS'
    $ws.Cells.Item(167, 12).Formula = '=ROUND((I167/$G$31-1)*100,3)+$L$29'

    # Row 168
    $ws.Cells.Item(168, 1).Value = 167
    $ws.Cells.Item(168, 2).Value = 150
    $ws.Cells.Item(168, 3).Value = 45264
    $ws.Cells.Item(168, 4).Value = '2023-12-04'
    $ws.Cells.Item(168, 5).Value = 1
    $ws.Cells.Item(168, 6).Value = 1.16
    $ws.Cells.Item(168, 7).Formula = '=I167'
    $ws.Cells.Item(168, 8).Value = 800
    $ws.Cells.Item(168, 9).Formula = '=G168+H168'
    $ws.Cells.Item(168, 10).Value = 'TENIS DE MESA'
    $ws.Cells.Item(168, 11).Value = 'MASTERS WOMEN'
    $ws.Cells.Item(168, 12).Formula = '=ROUND((I168/$G$31-1)*100,3)+$L$29'

    # Row 169
    $ws.Cells.Item(169, 1).Value = 168
    $ws.Cells.Item(169, 2).Value = 151
    $ws.Cells.Item(169, 3).Value = 45264
    $ws.Cells.Item(169, 4).Value = '2023-12-04'
    $ws.Cells.Item(169, 5).Value = 1
    $ws.Cells.Item(169, 6).Value = 1.165
    $ws.Cells.Item(169, 7).Formula = '=I168'
    $ws.Cells.Item(169, 8).Value = 297
    $ws.Cells.Item(169, 9).Formula = '=G169+H169'
    $ws.Cells.Item(169, 10).Value = 'TENIS DE MESA'
    $ws.Cells.Item(169, 11).Value = 'COPA TT POLONIA'
    $ws.Cells.Item(169, 12).Formula = '=ROUND((I169/$G$31-1)*100,3)+$L$29'

    # Row 170
    $ws.Cells.Item(170, 1).Value = 169
    $ws.Cells.Item(170, 2).Value = 152
    $ws.Cells.Item(170, 3).Value = 45264
    $ws.Cells.Item(170, 4).Value = '2023-12-04'
    $ws.Cells.Item(170, 5).Value = 0
    $ws.Cells.Item(170, 6).Value = 1.765
    $ws.Cells.Item(170, 7).Formula = '=I169'
    $ws.Cells.Item(170, 8).Value = -6097
    $ws.Cells.Item(170, 9).Formula = '=G170+H170'
    $ws.Cells.Item(170, 10).Value = 'ESPORTS'
    $ws.Cells.Item(170, 11).Value = 'LOL EUROPEAN CIRCUIT'
    $ws.Cells.Item(170, 12).Formula = '=ROUND((I170/$G$31-1)*100,3)+$L$29'

    # Row 171
    $ws.Cells.Item(171, 1).Value = 170
    $ws.Cells.Item(171, 2).Value = 153
    $ws.Cells.Item(171, 3).Value = 45264
    $ws.Cells.Item(171, 4).Value = '2023-12-04'
    $ws.Cells.Item(171, 5).Value = 1
    $ws.Cells.Item(171, 6).Value = 1.7
    $ws.Cells.Item(171, 7).Formula = '=I170'
    $ws.Cells.Item(171, 8).Value = 7000
    $ws.Cells.Item(171, 9).Formula = '=G171+H171'
    $ws.Cells.Item(171, 10).Value = 'ESPORTS'
    $ws.Cells.Item(171, 11).Value = 'LOL EUROPEAN CIRCUIT'
    $ws.Cells.Item(171, 12).Formula = '=ROUND((I171/$G$31-1)*100,3)+$L$29'


# Update the view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 152
$ws.Range("L173").Select()
